$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MissingWithItemForList")
$ws.Range("A11").Value = "Then"
